# Swap the presentation's applied design theme from "Integral" back to the
# default "Office Theme" palette.
#
# The deck's single slide master (and therefore every slide, since they all
# just inherit the master's color map) is tied to ppt/theme/theme1.xml. The
# font scheme and format scheme (fills/lines/effects) are identical between
# the "Integral" and "Office Theme" theme parts in this deck, so the only
# functional difference is the 12-colour `clrScheme`. Re-pointing those 12
# slots at the stock Office palette reproduces the authored edit.
#
# PowerPoint's VBA object model doesn't expose per-RGB setters directly on
# ColorScheme, but it does via the newer ThemeColorScheme collection
# (1=dk1, 2=lt1, 3=dk2, 4=lt2, 5-10=accent1-6, 11=hlink, 12=folHlink), which
# every Slide (and therefore the whole deck, since there's only one master)
# shares.

$p = $ppt.ActivePresentation

function HexToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Stock "Office Theme" colour scheme, in ThemeColorScheme index order.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$s = $p.Slides.Item(1)
for ($i = 1; $i -le 12; $i++) {
    $s.ThemeColorScheme.Item($i).RGB = HexToVbaRgb($officeColors[$i - 1])
}
